$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h volume-change snapshot (and the three
# reshuffled rows 41-43) to match the latest coinranking.com scrape.

# Row 2
$ws.Range("D2").Value = "'69.099.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.57%  "
# Row 3
$ws.Range("D3").Value = "'3.740.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.62%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'614.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
# Row 6
$ws.Range("D6").Value = "'177.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
# Row 7
$ws.Range("D7").Value = "'3.737.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
# Row 8
$ws.Range("E8").Value = "  +0.02%  "
# Row 9
$ws.Range("D9").Value = "'0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.00%  "
# Row 10
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.29%  "
# Row 11
$ws.Range("D11").Value = "'6.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "
# Row 12
$ws.Range("D12").Value = "'0.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.50%  "
# Row 13
$ws.Range("D13").Value = "'39.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.89%  "
# Row 14
$ws.Range("E14").Value = "  -1.91%  "
# Row 15
$ws.Range("D15").Value = "'4.364.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "
# Row 16
$ws.Range("D16").Value = "'3.742.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "
# Row 17
$ws.Range("D17").Value = "'69.207.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.46%  "
# Row 18
$ws.Range("E18").Value = "  -2.90%  "
# Row 19
$ws.Range("D19").Value = "'7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.25%  "
# Row 20
$ws.Range("D20").Value = "'16.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "
# Row 21
$ws.Range("D21").Value = "'497.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.33%  "
# Row 22
$ws.Range("D22").Value = "'9.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.04%  "
# Row 23
$ws.Range("E23").Value = "  -1.23%  "
# Row 24
$ws.Range("D24").Value = "'2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
# Row 25
$ws.Range("D25").Value = "'85.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.84%  "
# Row 26
$ws.Range("D26").Value = "'12.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.74%  "
# Row 27
$ws.Range("D27").Value = "'10.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.65%  "
# Row 28
$ws.Range("E28").Value = "  -4.59%  "
# Row 29
$ws.Range("E29").Value = "  +0.14%  "
# Row 30
$ws.Range("D30").Value = "'2.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
# Row 31
$ws.Range("E31").Value = "  +3.19%  "
# Row 32
$ws.Range("D32").Value = "'7.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.40%  "
# Row 33
$ws.Range("D33").Value = "'30.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.37%  "
# Row 34
$ws.Range("E34").Value = "  -1.96%  "
# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
# Row 36
$ws.Range("E36").Value = "  -0.56%  "
# Row 37
$ws.Range("D37").Value = "'6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "
# Row 38
$ws.Range("D38").Value = "'0.344"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "
# Row 39
$ws.Range("E39").Value = "  +4.30%  "
# Row 40
$ws.Range("D40").Value = "'451.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.70%  "
# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.35%  "
# Row 42
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'49.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.21%  "
# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.68%  "
# Row 44
$ws.Range("D44").Value = "'44.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
# Row 45
$ws.Range("D45").Value = "'8.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.88%  "
# Row 46
$ws.Range("D46").Value = "'2.939.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "
# Row 47
$ws.Range("D47").Value = "'0.0358"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
# Row 48
$ws.Range("D48").Value = "'27.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "
# Row 50
$ws.Range("D50").Value = "'137.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.63%  "
# Row 51
$ws.Range("E51").Value = "  -1.14%  "
